$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.785.69'
$ws.Range("E2").Value = '  -1.01%  '

$ws.Range("D3").Value = '2.341.23'
$ws.Range("E3").Value = '  -0.74%  '

$ws.Range("E4").Value = '  +0.22%  '

$ws.Range("E5").Value = '  -1.78%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '237.89'
$ws.Range("E6").Value = '  -0.41%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '72.52'
$ws.Range("E7").Value = '  -3.45%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  +6.57%  '

$ws.Range("E10").Value = '  -2.84%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.17'
$ws.Range("E11").Value = '  -0.38%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '31.95'
$ws.Range("E12").Value = '  +4.03%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.107'
$ws.Range("E13").Value = '  +0.43%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.18'
$ws.Range("E14").Value = '  +4.54%  '

$ws.Range("D15").Value = '2.689.51'
$ws.Range("E15").Value = '  -0.61%  '

$ws.Range("E16").Value = '  -4.34%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.888'
$ws.Range("E17").Value = '  -2.30%  '

$ws.Range("D18").Value = '2.330.70'
$ws.Range("E18").Value = '  -1.44%  '

$ws.Range("D19").Value = '43.579.88'
$ws.Range("E19").Value = '  -1.46%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.84'
$ws.Range("E20").Value = '  +4.71%  '

$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000100'
$ws.Range("E21").Value = '  -2.10%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '76.27'
$ws.Range("E22").Value = '  -2.68%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '255.50'
$ws.Range("E23").Value = '  -0.16%  '

$ws.Range("E24").Value = '  +22.50%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.01%  '

$ws.Range("E26").Value = '  -2.40%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.45'
$ws.Range("E27").Value = '  -2.86%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.55'
$ws.Range("E28").Value = '  +1.06%  '

$ws.Range("E29").Value = '  -0.82%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.48'
$ws.Range("E30").Value = '  -0.80%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '174.49'
$ws.Range("E31").Value = '  +0.44%  '

$ws.Range("E32").Value = '  -2.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.135'
$ws.Range("E33").Value = '  +1.37%  '

$ws.Range("E34").Value = '  +0.39%  '

$ws.Range("E35").Value = '  +6.77%  '

$ws.Range("E36").Value = '  -1.77%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.72'
$ws.Range("E37").Value = '  -3.76%  '

$ws.Range("E38").Value = '  -4.12%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.16'
$ws.Range("E39").Value = '  -5.09%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0275'
$ws.Range("E40").Value = '  +1.26%  '

$ws.Range("E41").Value = '  +10.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.203'
$ws.Range("E42").Value = '  +8.83%  '

$ws.Range("E43").Value = '  +0.01%  '

$ws.Range("B44").Value = 'MultiversX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '60.72'
$ws.Range("E44").Value = '  +15.07%  '

$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.90'
$ws.Range("E45").Value = '  +0.11%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.59'
$ws.Range("E46").Value = '  -3.65%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.66'
$ws.Range("E47").Value = '  +4.23%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.23'
$ws.Range("E48").Value = '  -2.67%  '

$ws.Range("E49").Value = '  +2.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '99.17'
$ws.Range("E50").Value = '  -0.74%  '

$ws.Range("E51").Value = '  -2.17%  '
